$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.034.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.84%  "

$ws.Range("D3").Value = "'2.969.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.42%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'594.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("D6").Value = "'149.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.04%  "

$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").Value = "'2.969.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.61%  "

$ws.Range("E9").Value = "  +1.37%  "

$ws.Range("E10").Value = "  +5.77%  "

$ws.Range("D11").Value = "'0.154"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.22%  "

$ws.Range("D12").Value = "'0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.30%  "

$ws.Range("D13").Value = "'0.0000242"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.55%  "

$ws.Range("D14").Value = "'33.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("E15").Value = "  -0.70%  "

$ws.Range("D16").Value = "'3.466.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.43%  "

$ws.Range("D17").Value = "'63.045.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.90%  "

$ws.Range("D18").Value = "'6.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("D19").Value = "'2.964.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.20%  "

$ws.Range("D20").Value = "'444.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.95%  "

$ws.Range("D21").Value = "'13.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.68%  "

$ws.Range("D22").Value = "'0.674"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.82%  "

$ws.Range("D23").Value = "'7.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("D24").Value = "'11.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.96%  "

$ws.Range("D25").Value = "'81.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").Value = "'2.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.11%  "

$ws.Range("D27").Value = "'11.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.29%  "

$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("E29").Value = "  +0.88%  "

$ws.Range("E30").Value = "  +21.62%  "

$ws.Range("D31").Value = "'7.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.43%  "

$ws.Range("D32").Value = "'2.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.17%  "

$ws.Range("D33").Value = "'26.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.45%  "

$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").Value = "'3.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.76%  "

$ws.Range("D37").Value = "'0.997"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.59%  "

$ws.Range("D38").Value = "'5.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.24%  "

$ws.Range("D39").Value = "'2.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.71%  "

$ws.Range("D40").Value = "'49.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("D41").Value = "'8.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.37%  "

$ws.Range("E42").Value = "  -4.28%  "

$ws.Range("D43").Value = "'0.285"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.93%  "

$ws.Range("D44").Value = "'40.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.44%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'372.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.05%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "'2.713.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.30%  "

$ws.Range("D47").Value = "'0.0342"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.55%  "

$ws.Range("D48").Value = "'135.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.55%  "

$ws.Range("D50").Value = "'23.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.24%  "

$ws.Range("D51").Value = "'0.105"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.28%  "
